$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowValues {
    param($Sheet, $Row, $Values)
    $n = $Values.Length
    $arr = New-Object 'object[,]' 1, $n
    for ($i = 0; $i -lt $n; $i++) {
        $arr[0, $i] = $Values[$i]
    }
    $startCell = $Sheet.Cells.Item($Row, 1)
    $endCell = $Sheet.Cells.Item($Row, $n)
    $range = $Sheet.Range($startCell, $endCell)
    $range.Value = $arr
}

# Row 2: sending cluster FAPs, Ccl21b -> Ccr10, target ECs (was Neutrophils); TPM values updated
Set-RowValues $ws 2 @(
    "FAPs", "Ccl21b", "Ccr10", "ECs",
    2,
    0.6666666666666666,
    0.126417,
    0.379251,
    0.733230478333749,
    0.8047954435010685,
    1,
    0.5,
    0.2651005,
    0.530201,
    0.2682690254597046,
    0.2682690254597046,
    0.0335132099085,
    0.201079259451,
    0.1967030258599479,
    0.2159016893224424
)

# Row 3: sending cluster FAPs (was MuSCs), Ccl21b -> Ccr10, target MuSCs (was Neutrophils); TPM values updated
Set-RowValues $ws 3 @(
    "FAPs", "Ccl21b", "Ccr10", "MuSCs",
    2,
    0.6666666666666666,
    0.126417,
    0.379251,
    0.733230478333749,
    0.8047954435010685,
    1,
    0.5,
    0.7230885,
    1.446177,
    0.7317309745402955,
    0.7317309745402955,
    0.0914106789045,
    0.548464073427,
    0.5365274524738012,
    0.5888937541786261
)

# Row 4 (new): sending cluster MuSCs, Ccl21b -> Ccr10, target ECs
Set-RowValues $ws 4 @(
    "MuSCs", "Ccl21b", "Ccr10", "ECs",
    1,
    0.5,
    0.045994,
    0.091988,
    0.266769521666251,
    0.1952045564989315,
    1,
    0.5,
    0.2651005,
    0.530201,
    0.2682690254597046,
    0.2682690254597046,
    0.012193032397,
    0.048772129588,
    0.07156599959975671,
    0.05236733613726221
)

# Row 5 (new): sending cluster MuSCs, Ccl21b -> Ccr10, target MuSCs
Set-RowValues $ws 5 @(
    "MuSCs", "Ccl21b", "Ccr10", "MuSCs",
    1,
    0.5,
    0.045994,
    0.091988,
    0.266769521666251,
    0.1952045564989315,
    1,
    0.5,
    0.7230885,
    1.446177,
    0.7317309745402955,
    0.7317309745402955,
    0.033257732469,
    0.133030929876,
    0.1952035220664943,
    0.1428372203616693
)

Write-Output ("Done. UsedRange: " + $ws.UsedRange.Address())
